$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - RandomForestClassifier
$ws.Range("C2").Value = 0.7796232301019086
$ws.Range("D2").Value = 0.7386778850941065
$ws.Range("E2").Value = 0.8190648727158267
$ws.Range("F2").Value = 0.9196525515743756
$ws.Range("G2").Value = 0.6395939086294417
$ws.Range("H2").Value = 0.9226579520697168
$ws.Range("I2").Value = 0.9196525515743756
$ws.Range("J2").Value = 0.921152800435019

# Row 3 - XGBClassifier
$ws.Range("C3").Value = 0.704335389143339
$ws.Range("D3").Value = 0.6660435680764363
$ws.Range("E3").Value = 0.7402250480598082
$ws.Range("F3").Value = 0.9771986970684039
$ws.Range("G3").Value = 0.4314720812182741
$ws.Range("H3").Value = 0.8893280632411067
$ws.Range("I3").Value = 0.9771986970684039
$ws.Range("J3").Value = 0.9311950336264874

# Row 4 - LogisticRegression
$ws.Range("C4").Value = 0.7658333195544459
$ws.Range("D4").Value = 0.7255832635751743
$ws.Range("E4").Value = 0.8041121902961073
$ws.Range("F4").Value = 0.9326818675352877
$ws.Range("G4").Value = 0.5989847715736041
$ws.Range("H4").Value = 0.9157782515991472
$ws.Range("I4").Value = 0.9326818675352877
$ws.Range("J4").Value = 0.9241527703066165
